$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.959.69'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.554.41'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').Value = '206.72'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').Value = '21.65'
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').Value = '0.0588'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').Value = '0.0859'
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = '1.776.39'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '1.554.87'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = '0.515'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = '26.953.63'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '61.80'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').Value = '214.63'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '7.25'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('D23').Value = '9.19'
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('D25').Value = '153.31'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '1.09'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('E32').Value = '  +1.53%  '
$ws.Range('D33').Value = '1.382.80'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('D36').Value = '0.972'
$ws.Range('E36').Value = '  +5.69%  '
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('E38').Value = '  +1.83%  '
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('E41').Value = '  +0.52%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '2.26'
$ws.Range('E43').Value = '  +3.10%  '
$ws.Range('D44').Value = '5.47'
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('D45').Value = '63.90'
$ws.Range('E45').Value = '  +1.44%  '
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').Value = '1.690.11'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = '86.04'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').Value = '0.0957'
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₇0959'
$ws.Range('E51').Value = '  -1.29%  '
